$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = [double]"25.00000000000047"
$ws.Range("H2").Value = [double]"2.163442525127124e-08"
$ws.Range("I2").Value = [double]"2.163442525127124e-08"
$ws.Range("L2").Value = [double]"40.59733642051505"
$ws.Range("M2").Value = "[28.39046752180237, 52.80420531922773]"
$ws.Range("N2").Value = [double]"2.857147607393529e-08"
$ws.Range("O2").Value = [double]"2.857147607393529e-08"
$ws.Range("P2").Value = [double]"1.553500271144503"
$ws.Range("Q2").Value = "[1.1887107337907326, 1.9182898084982734]"
$ws.Range("R2").Value = [double]"5.057199103930543e-11"
$ws.Range("S2").Value = [double]"5.057199103930543e-11"
$ws.Range("T2").Value = [double]"45.67476488200646"
$ws.Range("U2").Value = "[37.62321548104497, 53.72631428296794]"
$ws.Range("V2").Value = [double]"6.883382752675971e-15"
$ws.Range("W2").Value = [double]"6.883382752675971e-15"
$ws.Range("X2").Value = [double]"18.81881881881917"
$ws.Range("Y2").Value = [double]"17.36736736736768"
$ws.Range("Z2").Value = [double]"20.27027027027065"
$ws.Range("F3").Value = [double]"25.00000000000047"
$ws.Range("H3").Value = [double]"8.201268220098257e-09"
$ws.Range("I3").Value = [double]"8.201268220098257e-09"
$ws.Range("L3").Value = [double]"48.09179624470279"
$ws.Range("M3").Value = "[32.23993616846488, 63.943656320940704]"
$ws.Range("N3").Value = [double]"2.140216388557548e-07"
$ws.Range("O3").Value = [double]"2.140216388557548e-07"
$ws.Range("P3").Value = [double]"1.314500229429964"
$ws.Range("Q3").Value = "[0.9622896416401163, 1.666710817219811]"
$ws.Range("R3").Value = [double]"1.755658285418349e-09"
$ws.Range("S3").Value = [double]"1.755658285418349e-09"
$ws.Range("T3").Value = [double]"54.39633553393848"
$ws.Range("U3").Value = "[45.274612912458366, 63.5180581554186]"
$ws.Range("V3").Value = [double]"1.110223024625157e-15"
$ws.Range("W3").Value = [double]"1.110223024625157e-15"
$ws.Range("X3").Value = [double]"19.76976976977014"
$ws.Range("Y3").Value = [double]"18.36836836836871"
$ws.Range("Z3").Value = [double]"21.17117117117157"
$ws.Range("F4").Value = [double]"25.00000000000047"
$ws.Range("H4").Value = [double]"2.252333474173085e-07"
$ws.Range("I4").Value = [double]"2.252333474173085e-07"
$ws.Range("L4").Value = [double]"39.46604517489112"
$ws.Range("M4").Value = "[24.489831547125156, 54.442258802657086]"
$ws.Range("N4").Value = [double]"3.28159647033921e-06"
$ws.Range("O4").Value = [double]"3.28159647033921e-06"
$ws.Range("P4").Value = [double]"1.352237078121732"
$ws.Range("Q4").Value = "[0.9119738433844242, 1.7925003128590395]"
$ws.Range("R4").Value = [double]"1.651552159476211e-07"
$ws.Range("S4").Value = [double]"1.651552159476211e-07"
$ws.Range("T4").Value = [double]"52.87604776505727"
$ws.Range("U4").Value = "[44.2591734032096, 61.49292212690493]"
$ws.Range("V4").Value = [double]"4.440892098500626e-16"
$ws.Range("W4").Value = [double]"4.440892098500626e-16"
$ws.Range("X4").Value = [double]"19.61961961961999"
$ws.Range("Y4").Value = [double]"17.86786786786821"
$ws.Range("Z4").Value = [double]"21.37137137137178"
$ws.Range("F5").Value = [double]"25.00000000000047"
$ws.Range("H5").Value = [double]"5.45724887046628e-06"
$ws.Range("I5").Value = [double]"5.45724887046628e-06"
$ws.Range("L5").Value = [double]"37.51186785259279"
$ws.Range("M5").Value = "[19.092245969111552, 55.931489736074035]"
$ws.Range("N5").Value = [double]"0.0001697083080771655"
$ws.Range("O5").Value = [double]"0.0001697083080771655"
$ws.Range("P5").Value = [double]"0.7987632973091161"
$ws.Range("Q5").Value = "[0.34592111300788453, 1.2516054816103477]"
$ws.Range("R5").Value = [double]"0.000908661047379411"
$ws.Range("S5").Value = [double]"0.000908661047379411"
$ws.Range("T5").Value = [double]"48.8985040676629"
$ws.Range("U5").Value = "[39.53126397628262, 58.26574415904318]"
$ws.Range("V5").Value = [double]"1.056932319443149e-13"
$ws.Range("W5").Value = [double]"1.056932319443149e-13"
$ws.Range("X5").Value = [double]"21.82182182182223"
$ws.Range("Y5").Value = [double]"20.0200200200204"
$ws.Range("Z5").Value = [double]"23.62362362362407"
$ws.Range("F6").Value = [double]"24.33000000000036"
$ws.Range("H6").Value = [double]"3.963274153306884e-12"
$ws.Range("I6").Value = [double]"3.963274153306884e-12"
$ws.Range("L6").Value = [double]"51.9166441207228"
$ws.Range("M6").Value = "[38.55920755467403, 65.27408068677157]"
$ws.Range("N6").Value = [double]"6.136280372714964e-10"
$ws.Range("O6").Value = [double]"6.136280372714964e-10"
$ws.Range("P6").Value = [double]"0.3585000625718084"
$ws.Range("Q6").Value = "[0.09434212172942402, 0.6226580034141929]"
$ws.Range("R6").Value = [double]"0.008928589696802103"
$ws.Range("S6").Value = [double]"0.008928589696802103"
$ws.Range("T6").Value = [double]"57.10809535188693"
$ws.Range("U6").Value = "[49.746317526775215, 64.46987317699865]"
$ws.Range("X6").Value = [double]"22.94180180180214"
$ws.Range("Y6").Value = [double]"21.91891891891925"
$ws.Range("Z6").Value = [double]"23.96468468468504"
$ws.Range("F7").Value = [double]"24.33000000000036"
$ws.Range("H7").Value = [double]"2.681529244819014e-05"
$ws.Range("I7").Value = [double]"2.681529244819014e-05"
$ws.Range("L7").Value = [double]"39.79674538217294"
$ws.Range("M7").Value = "[20.452863361151692, 59.14062740319418]"
$ws.Range("N7").Value = [double]"0.0001487367425228037"
$ws.Range("O7").Value = [double]"0.0001487367425228037"
$ws.Range("P7").Value = [double]"0.4213948103914236"
$ws.Range("Q7").Value = "[-0.11950002085726918, 0.9622896416401163]"
$ws.Range("R7").Value = [double]"0.1236241087023064"
$ws.Range("S7").Value = [double]"0.1236241087023064"
$ws.Range("T7").Value = [double]"47.35415561173416"
$ws.Range("U7").Value = "[36.4807208135709, 58.22759040989743]"
$ws.Range("V7").Value = [double]"2.671440846313544e-11"
$ws.Range("W7").Value = [double]"2.671440846313544e-11"
$ws.Range("X7").Value = [double]"22.6982582582586"
$ws.Range("Y7").Value = [double]"20.60378378378409"
$ws.Range("Z7").Value = [double]"24.7927327327331"
$ws.Range("F8").Value = [double]"24.33000000000036"
$ws.Range("H8").Value = [double]"2.343109484215233e-09"
$ws.Range("I8").Value = [double]"2.343109484215233e-09"
$ws.Range("L8").Value = [double]"48.79360900885217"
$ws.Range("M8").Value = "[34.158801081895746, 63.428416935808585]"
$ws.Range("N8").Value = [double]"2.69822741927328e-08"
$ws.Range("O8").Value = [double]"2.69822741927328e-08"
$ws.Range("P8").Value = [double]"0.1823947686768852"
$ws.Range("Q8").Value = "[-0.1320789704211922, 0.4968685077749626]"
$ws.Range("R8").Value = [double]"0.2488860963610899"
$ws.Range("S8").Value = [double]"0.2488860963610899"
$ws.Range("T8").Value = [double]"58.50267738080046"
$ws.Range("U8").Value = "[49.86705973898089, 67.13829502262004]"
$ws.Range("X8").Value = [double]"23.62372372372408"
$ws.Range("Y8").Value = [double]"22.40600600600634"
$ws.Range("Z8").Value = [double]"24.84144144144181"
$ws.Range("F9").Value = [double]"24.33000000000036"
$ws.Range("H9").Value = [double]"3.237332635297463e-07"
$ws.Range("I9").Value = [double]"3.237332635297463e-07"
$ws.Range("L9").Value = [double]"41.30772619293448"
$ws.Range("M9").Value = "[26.04179608457298, 56.57365630129598]"
$ws.Range("N9").Value = [double]"2.030020677556976e-06"
$ws.Range("O9").Value = [double]"2.030020677556976e-06"
$ws.Range("P9").Value = [double]"0.1823947686768852"
$ws.Range("Q9").Value = "[-0.25786846606042246, 0.6226580034141929]"
$ws.Range("R9").Value = [double]"0.40845545369119"
$ws.Range("S9").Value = [double]"0.40845545369119"
$ws.Range("T9").Value = [double]"56.31644551507333"
$ws.Range("U9").Value = "[47.34747801621035, 65.2854130139363]"
$ws.Range("V9").Value = [double]"2.220446049250313e-16"
$ws.Range("W9").Value = [double]"2.220446049250313e-16"
$ws.Range("X9").Value = [double]"23.62372372372408"
$ws.Range("Y9").Value = [double]"21.91891891891925"
$ws.Range("Z9").Value = [double]"25.3285285285289"
